$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order row (row 10): Order ID, Order Date, Order Status.
# Force text storage (matching the existing text-stored ID/date columns)
# by temporarily applying a text number format, then restore the default
# "Normal" style so the cells end up unstyled like the rest of the data rows.
$ws.Range("A10:C10").NumberFormat = "@"
$ws.Range("A10").Value = "1920961"
$ws.Range("B10").Value = "03/11/2025"
$ws.Range("C10").Value = "Pending"
$ws.Range("A10:C10").Style = "Normal"
